$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category") between the existing
# "total" (G) and "date" (old H, now shifts to I) columns. This pushes
# date / legislator_name / legislator_id from H/I/J to I/J/K.
$ws.Columns.Item(8).Insert()

# Header for the new column.
$ws.Range("H1").Value = "property_category"

# Every stock record in this sheet is a "stock" property, so fill the
# new column for all 18 data rows.
$ws.Range("H2:H19").Value = "stock"

# Data cleanup: the "total" value for row 17 was stored as text with a
# full-width comma ("45，320"); normalize it to a plain "45320" string
# while keeping it a text cell (not a number) and preserving its
# existing cell style.
$ws.Range("G17").Formula = "=""45320"""
$ws.Range("G17").Copy()
$ws.Range("G17").PasteSpecial(-4163)
$excel.CutCopyMode = $false
